$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the existing 8-row block (155:162) down to the new block (163:170)
# so the new rows inherit identical cell styling (alignment / red "stale" font)
# without registering any new style entries.
$ws.Range("A155:H162").Copy()
$ws.Range("A163:H170").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 163
$ws.Range("A163").Value = "2026-02-22 17:49:12"
$ws.Range("B163").Value = "poqui"
$ws.Range("C163").Value = "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza"
$ws.Range("D163").Value = 2049
$ws.Range("E163").Value = ""
$ws.Range("F163").Value = ""
$ws.Range("G163").Value = "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html"
$ws.Range("H163").Value = "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc"

# Row 164
$ws.Range("A164").Value = "2026-02-22 17:49:12"
$ws.Range("B164").Value = "poqui"
$ws.Range("C164").Value = "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda"
$ws.Range("D164").Value = 2299
$ws.Range("E164").Value = "19.01.2026"
$ws.Range("F164").Value = 34
$ws.Range("G164").Value = "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html"
$ws.Range("H164").Value = "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR"

# Row 165
$ws.Range("A165").Value = "2026-02-22 17:49:12"
$ws.Range("B165").Value = "poqui"
$ws.Range("C165").Value = "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy"
$ws.Range("D165").Value = 2499
$ws.Range("E165").Value = "28.10.2025"
$ws.Range("F165").Value = 117
$ws.Range("G165").Value = "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html"
$ws.Range("H165").Value = "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger"

# Row 166
$ws.Range("A166").Value = "2026-02-22 17:49:12"
$ws.Range("B166").Value = "poqui"
$ws.Range("C166").Value = "Przytulny pokój blisko Politechniki – ul. Przytulna"
$ws.Range("D166").Value = 549
$ws.Range("E166").Value = "10.10.2025"
$ws.Range("F166").Value = 135
$ws.Range("G166").Value = "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html"
$ws.Range("H166").Value = "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz"

# Row 167
$ws.Range("A167").Value = "2026-02-22 17:49:12"
$ws.Range("B167").Value = "pokojewlublinie"
$ws.Range("C167").Value = "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58"
$ws.Range("D167").Value = 0
$ws.Range("E167").Value = "11.08.2025"
$ws.Range("F167").Value = 195
$ws.Range("G167").Value = "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html"
$ws.Range("H167").Value = "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm"

# Row 168
$ws.Range("A168").Value = "2026-02-22 17:49:12"
$ws.Range("B168").Value = "pokojewlublinie"
$ws.Range("C168").Value = "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12"
$ws.Range("D168").Value = 12640
$ws.Range("E168").Value = "19.01.2026"
$ws.Range("F168").Value = 34
$ws.Range("G168").Value = "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html"
$ws.Range("H168").Value = "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc"

# Row 169
$ws.Range("A169").Value = "2026-02-22 17:49:12"
$ws.Range("B169").Value = "dawnypatron"
$ws.Range("C169").Value = "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4."
$ws.Range("D169").Value = 730
$ws.Range("E169").Value = "20.09.2024"
$ws.Range("F169").Value = 520
$ws.Range("G169").Value = "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html"
$ws.Range("H169").Value = "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM"

# Row 170
$ws.Range("A170").Value = "2026-02-22 17:49:12"
$ws.Range("B170").Value = "dawnypatron"
$ws.Range("C170").Value = "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14"
$ws.Range("D170").Value = 14690
$ws.Range("E170").Value = "05.12.2025"
$ws.Range("F170").Value = 79
$ws.Range("G170").Value = "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html"
$ws.Range("H170").Value = "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv"
